$d = $word.ActiveDocument

function Replace-Text($findText, $replaceText) {
    $found = $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false,
                                      $true, 1, $false, $replaceText, 2)
    if (-not $found) {
        throw "Find/Replace failed to locate: $findText"
    }
}

# ---------------------------------------------------------------------------
# 1. Remove the _GoBack bookmark from its original location (between "does"
#    and " not work." in the title paragraph). Word re-anchors this bookmark
#    to the site of the most recent edit, which we recreate in step 6 below.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 2. "Setting up ORNG in Profiles is complex." -> "...Profiles RNS is complex."
# ---------------------------------------------------------------------------
Replace-Text "Setting up ORNG in Profiles is complex" "Setting up ORNG in Profiles RNS is complex"

# ---------------------------------------------------------------------------
# 3a. "...make sense for your Profiles server." -> "...your Profiles RNS server."
# ---------------------------------------------------------------------------
Replace-Text "make sense for your Profiles server" "make sense for your Profiles RNS server"

# ---------------------------------------------------------------------------
# 3b. "Sometimes your profiles server will" -> "Sometimes your server will"
# ---------------------------------------------------------------------------
Replace-Text "Sometimes your profiles server will" "Sometimes your server will"

# ---------------------------------------------------------------------------
# 4a. "...an ORNG application on a Profiles web page," -> "...Profiles RNS web page,"
# ---------------------------------------------------------------------------
Replace-Text "application on a Profiles web page" "application on a Profiles RNS web page"

# ---------------------------------------------------------------------------
# 4b. "it's the Profiles web server." -> "it's the Profiles RNS web server."
# ---------------------------------------------------------------------------
$rightQuote = [char]0x2019
Replace-Text ($rightQuote + "s the Profiles web server") ($rightQuote + "s the Profiles RNS web server")

# ---------------------------------------------------------------------------
# 4c. "...example ones provided with Profiles are free..." -> "...Profiles RNS are free..."
# ---------------------------------------------------------------------------
Replace-Text "provided with Profiles are free" "provided with Profiles RNS are free"

# ---------------------------------------------------------------------------
# 5. "Use the sandbox at http://[Your Profiles URL]/ORNG" -> "...Profiles RNS URL]/ORNG"
# ---------------------------------------------------------------------------
Replace-Text "[Your Profiles URL]" "[Your Profiles RNS URL]"

# ---------------------------------------------------------------------------
# 6. Re-create the _GoBack bookmark around the word "Profiles" in the
#    sentence we just edited - this mirrors where real Word leaves the
#    hidden _GoBack bookmark: the site of the most recent edit.
# ---------------------------------------------------------------------------
$srch = $d.Content
$found = $srch.Find.Execute("[Your Profiles RNS URL]", $true, $false, $false, $false, $false,
                             $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not relocate '[Your Profiles RNS URL]' to place _GoBack bookmark"
}

$bmStart = $srch.Start + ("[Your ").Length
$bmEnd = $bmStart + ("Profiles").Length
$bmRange = $d.Range($bmStart, $bmEnd)
if ($bmRange.Text -ne "Profiles") {
    throw "Bookmark range text mismatch: $($bmRange.Text)"
}
$d.Bookmarks.Add("_GoBack", $bmRange)
